$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Child")

$updates = @{
    2  = "-6.65,7.8"
    3  = "7.75,-1.13"
    4  = "-8.38,6.18"
    5  = "-5.15,-3.61"
    6  = "-7.14,-1.26"
    7  = "-0.63,-1.53"
    8  = "6.04,5.4"
    9  = "5.19,-5.63"
    10 = "3.52,-5.3"
    11 = "-6.03,3.13"
    12 = "-6.3,-0.62"
    13 = "-9.16,-4.56"
    14 = "-0.18,-7.37"
    15 = "4.45,1.94"
    16 = "-7.85,2.56"
    17 = "-9.9,3.16"
    18 = "-4.83,-1.02"
    19 = "-6.72,-1.52"
    20 = "0.42,9.67"
    21 = "-5.19,6.9"
    22 = "-2.33,-7.05"
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
